# Apply the "Add files via upload" revision to the Global Advisor tracker
# slide:
#   1) TextBox 2 ("Translation specs...") - reword the second sentence and
#      merge the two paragraphs into one (shape autofits to a shorter
#      height as a result).
#   2) TextBox 11 ("Press release: ...") - rename the label to
#      "Public visibility: " and recolor it to the muted gray used
#      elsewhere in the deck.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------------
# 1) "TextBox 2" - Translation specs / Turnaround time paragraph
# ---------------------------------------------------------------------
$specsShape = $s.Shapes.Item("TextBox 2")
$specsRange = $specsShape.TextFrame.TextRange

$run1 = "Translation specs: 3,000 words/month, around 32 countries and 35 languages, human "
$run2 = "translation Turnaround "
$run3 = "time from handoff of final materials: 4 days"

# Replacing the whole text range collapses the two paragraphs into a
# single one (the shape is spAutoFit, so its height recalculates too).
$specsRange.Text = $run1 + $run2 + $run3

# Re-split into three runs by nudging the formatting per character range
# so the OOXML keeps the same run boundaries as the authored edit.
$specsShape.TextFrame.TextRange.Characters(1, $run1.Length).Font.Color.RGB = 0x3E3E3E
$specsShape.TextFrame.TextRange.Characters($run1.Length + 1, $run2.Length).Font.Color.RGB = 0x3E3E3E
$specsShape.TextFrame.TextRange.Characters($run1.Length + $run2.Length + 1, $run3.Length).Font.Color.RGB = 0x3E3E3E

# ---------------------------------------------------------------------
# 2) "TextBox 11" - Press release label -> Public visibility label
# ---------------------------------------------------------------------
$visShape = $s.Shapes.Item("TextBox 11")
$visRange = $visShape.TextFrame.TextRange

$oldLabel = "Press release: "
$newLabel = "Public visibility: "

$visRange.Characters(1, $oldLabel.Length).Text = $newLabel
$visShape.TextFrame.TextRange.Characters(1, $newLabel.Length).Font.Color.RGB = 0x888888
